# Add support for string instruments: violin, viola, cello, double bass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (internal keys) -------------------------------------------
$ws.Range("A15").Value = "violin"
$ws.Range("A16").Value = "viola"
$ws.Range("A17").Value = "cello"
$ws.Range("A18").Value = "double_bass"

# --- Columns B:E (display names, same across en/it/de/lv for now) -------
$ws.Range("B15:E15").Value = "Violin"
$ws.Range("B16:E16").Value = "Viola"
$ws.Range("B17:E17").Value = "Cello"
$ws.Range("B18:E18").Value = "Double Bass"

# --- Columns F:H (low_note, high_note, transpose) ------------------------
$ws.Range("F15").Value = 55
$ws.Range("G15").Value = 93
$ws.Range("H15").Value = 0

$ws.Range("F16").Value = 48
$ws.Range("G16").Value = 84
$ws.Range("H16").Value = 0

$ws.Range("F17").Value = 36
$ws.Range("G17").Value = 84
$ws.Range("H17").Value = 0

$ws.Range("F18").Value = 28
$ws.Range("G18").Value = 67
$ws.Range("H18").Value = 0

# --- Column I (clef) ------------------------------------------------------
$ws.Range("I15").Value = "treble"
$ws.Range("I16").Value = "alto"
$ws.Range("I17").Value = "bass"
$ws.Range("I18").Value = "bass"

# Move the active selection, matching the saved cursor position.
$ws.Range("F30").Select() | Out-Null
